$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E3: 5 -> 1
$ws.Range("E3").Value = 1

# Row 5: update date/time (keep C5, D5, E5 unchanged)
$ws.Range("A5").Value = (Get-Date -Year 2021 -Month 1 -Day 25 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B5").Value = 0.45043981481481

# Row 6: update date/time, add shared string "Cargue Inicial" in C6, update D6
$ws.Range("A6").Value = (Get-Date -Year 2021 -Month 1 -Day 25 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B6").Value = 0.44269675925926
$ws.Range("C6").Value = "Cargue Inicial"
$ws.Range("D6").Value = 2
